$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Remove the two sample award entries (rows 2 and 3) from the Data sheet.
# Clearing (not just clearing contents) removes both value and formatting so
# the now-empty rows drop out of the saved sheetData entirely, while rows
# below (e.g. row 10) keep their original row numbers.
$ws.Range("A2:F3").Clear()

# Reflect where the user ended up clicking after removing the sample rows.
$ws.Activate()
[void]$ws.Range("C10").Select()
